{"js": "// Update the multiplication problems in the practice-sheet table.\n// Each \"from\" value occurs exactly once in the document, so a plain\n// text search locates the correct cell. All searches are performed\n// (and resolved) BEFORE any text is written back, because a couple of\n// the replacement values happen to equal another problem's original\n// text (e.g. \"615\u00d73=\" is both replaced away and later re-introduced as\n// a replacement elsewhere) - resolving the searches first avoids any\n// accidental re-matching of newly written text.\nconst mapping = [\n  { from: \"948\u00d72=\", to: \"544\u00d76=\" },\n  { from: \"769\u00d78=\", to: \"505\u00d76=\" },\n  { from: \"332\u00d77=\", to: \"320\u00d73=\" },\n  { from: \"928\u00d77=\", to: \"486\u00d74=\" },\n  { from: \"353\u00d75=\", to: \"774\u00d74=\" },\n  { from: \"355\u00d75=\", to: \"310\u00d76=\" },\n  { from: \"558\u00d78=\", to: \"291\u00d75=\" },\n  { from: \"582\u00d75=\", to: \"707\u00d76=\" },\n  { from: \"351\u00d74=\", to: \"617\u00d73=\" },\n  { from: \"615\u00d73=\", to: \"942\u00d72=\" },\n  { from: \"302\u00d73=\", to: \"316\u00d76=\" },\n  { from: \"775\u00d73=\", to: \"107\u00d75=\" },\n  { from: \"938\u00d76=\", to: \"829\u00d78=\" },\n  { from: \"133\u00d79=\", to: \"254\u00d76=\" },\n  { from: \"663\u00d78=\", to: \"639\u00d75=\" },\n  { from: \"195\u00d79=\", to: \"176\u00d72=\" },\n  { from: \"165\u00d72=\", to: \"710\u00d72=\" },\n  { from: \"461\u00d78=\", to: \"375\u00d75=\" },\n  { from: \"108\u00d78=\", to: \"615\u00d73=\" },\n  { from: \"456\u00d77=\", to: \"861\u00d79=\" },\n  { from: \"207\u00d79=\", to: \"850\u00d77=\" },\n  { from: \"619\u00d77=\", to: \"488\u00d77=\" },\n  { from: \"392\u00d72=\", to: \"350\u00d74=\" },\n  { from: \"692\u00d76=\", to: \"513\u00d79=\" },\n  { from: \"317\u00d78=\", to: \"494\u00d72=\" },\n];\n\n// Phase 1: search for every \"from\" string and load the results.\nconst searches = mapping.map((m) =>\n  context.document.body.search(m.from, { matchCase: true })\n);\nsearches.forEach((s) => s.load(\"items\"));\nawait context.sync();\n\n// Phase 2: replace each match with its corresponding \"to\" text, now\n// that all of the original text has already been located.\nfor (let i = 0; i < mapping.length; i++) {\n  const items = searches[i].items;\n  if (items.length !== 1) {\n    throw new Error(\n      `expected exactly one match for \"${mapping[i].from}\", found ${items.length}`\n    );\n  }\n  items[0].insertText(mapping[i].to, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"From\" value occurs exactly once in the document. A couple of\n# the \"To\" values happen to equal another problem's original \"From\"\n# text (e.g. \"615\u00d73=\" is replaced away at one cell and re-introduced as\n# a replacement at another), so this is done in two phases: first every\n# \"From\" string is located (recording its Start/End), and only once all\n# of the original locations are known are the replacement texts written\n# back. That way a freshly written value is never re-matched by a later\n# search. (Parallel arrays are used instead of mutating hashtable\n# entries from inside a foreach loop, since such in-loop mutations do\n# not reliably persist here.)\n$froms = @(\n    \"948\u00d72=\", \"769\u00d78=\", \"332\u00d77=\", \"928\u00d77=\", \"353\u00d75=\",\n    \"355\u00d75=\", \"558\u00d78=\", \"582\u00d75=\", \"351\u00d74=\", \"615\u00d73=\",\n    \"302\u00d73=\", \"775\u00d73=\", \"938\u00d76=\", \"133\u00d79=\", \"663\u00d78=\",\n    \"195\u00d79=\", \"165\u00d72=\", \"461\u00d78=\", \"108\u00d78=\", \"456\u00d77=\",\n    \"207\u00d79=\", \"619\u00d77=\", \"392\u00d72=\", \"692\u00d76=\", \"317\u00d78=\"\n)\n$tos = @(\n    \"544\u00d76=\", \"505\u00d76=\", \"320\u00d73=\", \"486\u00d74=\", \"774\u00d74=\",\n    \"310\u00d76=\", \"291\u00d75=\", \"707\u00d76=\", \"617\u00d73=\", \"942\u00d72=\",\n    \"316\u00d76=\", \"107\u00d75=\", \"829\u00d78=\", \"254\u00d76=\", \"639\u00d75=\",\n    \"176\u00d72=\", \"710\u00d72=\", \"375\u00d75=\", \"615\u00d73=\", \"861\u00d79=\",\n    \"850\u00d77=\", \"488\u00d77=\", \"350\u00d74=\", \"513\u00d79=\", \"494\u00d72=\"\n)\n\n$starts = @()\n$ends = @()\n\n# Phase 1: locate every \"From\" string.\nfor ($i = 0; $i -lt $froms.Count; $i++) {\n    $r = $d.Range()\n    $find = $r.Find\n    $find.Text = $froms[$i]\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"could not find text: $($froms[$i])\"\n    }\n    $starts += $r.Start\n    $ends += $r.End\n}\n\n# Phase 2: write the replacement text into each recorded location.\nfor ($i = 0; $i -lt $froms.Count; $i++) {\n    $target = $d.Range($starts[$i], $ends[$i])\n    $target.Text = $tos[$i]\n}\n"}
